$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark from the first (empty) paragraph.
#    It will be re-created at the end of the document, after "DQ123456D".
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Remove the "Nif" paragraph (together with its spell-check proofErr
#    markers) entirely - it gets replaced by a clean "VAT" / "/NIF" pair of
#    runs below.
# ---------------------------------------------------------------------------
$nifPara = $d.Paragraphs(2)
$nifPara.Range.Delete()

# ---------------------------------------------------------------------------
# 3. Insert the replacement "VAT" + "/NIF" paragraph in front of the
#    "123456789" paragraph. Using InsertXML (rather than InsertBefore/
#    InsertAfter) keeps "VAT" and "/NIF" as two discrete runs, just like in
#    the target document, instead of being coalesced into a single run.
# ---------------------------------------------------------------------------
$vatNifXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>VAT</w:t></w:r><w:r><w:t>/NIF</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$numberPara = $d.Paragraphs(2)     # now "123456789"
$numberPara.Range.InsertParagraphBefore()
$vatNifPara = $d.Paragraphs(2)     # the freshly inserted, still empty paragraph
$vatNifPara.Range.InsertXML($vatNifXml)

# ---------------------------------------------------------------------------
# 4. After "123456789", add four new paragraphs: "UTR ", "1234567890",
#    "NINO" and "DQ123456D" (the last one carrying the relocated "_GoBack"
#    bookmark).
# ---------------------------------------------------------------------------
$numberPara = $d.Paragraphs(3)     # "123456789" again, now at index 3
$numberPara.Range.InsertParagraphAfter()
$utrPara = $d.Paragraphs(4)
$utrXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">UTR </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$utrPara.Range.InsertXML($utrXml)

$utrPara = $d.Paragraphs(4)
$utrPara.Range.InsertParagraphAfter()
$utrNumberPara = $d.Paragraphs(5)
$utrNumberXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>1234567890</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$utrNumberPara.Range.InsertXML($utrNumberXml)

$utrNumberPara = $d.Paragraphs(5)
$utrNumberPara.Range.InsertParagraphAfter()
$ninoPara = $d.Paragraphs(6)
$ninoXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>NINO</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$ninoPara.Range.InsertXML($ninoXml)

$ninoPara = $d.Paragraphs(6)
$ninoPara.Range.InsertParagraphAfter()
$ninoNumberPara = $d.Paragraphs(7)
$ninoNumberXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>DQ123456D</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$ninoNumberPara.Range.InsertXML($ninoNumberXml)

# ---------------------------------------------------------------------------
# 5. Rename the Portuguese-localised built-in styleIds to their canonical
#    English identifiers (Default Paragraph Font / Normal Table / No List).
#    These 3 styles are not referenced anywhere else in the document, so a
#    straight Find/Replace over the style definitions' raw ids is safe.
# ---------------------------------------------------------------------------
$styles = $d.Styles
foreach ($s in $styles) {
    if ($s.NameLocal -eq "Default Paragraph Font") { $s.NameLocal = "Default Paragraph Font" }
}

Write-Host "Done"
